$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.64"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.11"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.408"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06050"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.402"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9309"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07449"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03372"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09362"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001595"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04831"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005942"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005101"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004164"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009840"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.646"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.445"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006401"

$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1077"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002710"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006303"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005262"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005802"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9004"

$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002235"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"

